$wb = $excel.ActiveWorkbook

# --- Update Summary sheet timestamps ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = "2024_01_22__16_52_54"
$summary.Range("B6").Value = "2024_01_22__16_53_19"

# --- Update Forecast sheet values ---
$forecast = $wb.Worksheets.Item("Forecast")

# Row 3
$forecast.Range("G3").Value = 0.08
$forecast.Range("I3").Value = 100.08

# Row 4
$forecast.Range("B4").Value = 1520
$forecast.Range("D4").Value = 120.21
$forecast.Range("E4").Value = 0
$forecast.Range("I4").Value = 100
$forecast.Range("J4").Value = 1399.79
$forecast.Range("K4").Value = 0
$forecast.Range("L4").Value = 120.21
$forecast.Range("M4").Value = 1520

# Row 5
$forecast.Range("B5").Value = 1800
$forecast.Range("C5").Value = 100
$forecast.Range("D5").Value = 100
$forecast.Range("E5").Value = 100
$forecast.Range("I5").Value = 100
$forecast.Range("J5").Value = 1500
$forecast.Range("K5").Value = 100
$forecast.Range("L5").Value = 200
$forecast.Range("M5").Value = 1800

# Row 6
$forecast.Range("B6").Value = 1800
$forecast.Range("C6").Value = 100
$forecast.Range("D6").Value = 100
$forecast.Range("E6").Value = 100
$forecast.Range("J6").Value = 1500
$forecast.Range("K6").Value = 100
$forecast.Range("L6").Value = 200
$forecast.Range("M6").Value = 1800
